# Sprint Backlog / Burndown workbook update
# 1) Restructured test packages (inserted new rows into the task list)
# 2) Fixed classpath / assignee values and filled in "Week 2" actuals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure rows: insert the new task rows so the table grows from
#     28 data rows to 31 data rows before the totals line ---

# New task inserted into the "Create Inventory" section (becomes row 21)
$ws.Rows("21:21").Insert()

# Two new tasks inserted into the "Create Orders" section (become rows 24:25
# after the previous insert already shifted things down by one)
$ws.Rows("24:25").Insert()

# --- Fix assignee for the "Product" button tasks: TBD -> Vitor ---
$ws.Range("B16:B20").Value = "Vitor"

# --- Seed the brand-new task descriptions first so they land in the shared
#     string table in the same order they were authored (replacing "TBD") ---
$ws.Range("C23").Value = "Implement functionality for orders"
$ws.Range("C24").Value = "Testing for orders functionality"
$ws.Range("C21").Value = "Create Model Class for inventory"

# --- Row by row content / numbers (Week columns) ---

$ws.Range("A3").Value = "Login"
$ws.Range("B3").Value = "Jason"
$ws.Range("C3").Value = "Create Login Page"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0

$ws.Range("A4").Value = "Login"
$ws.Range("B4").Value = "Jason"
$ws.Range("C4").Value = "Create credentials data base"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 1

$ws.Range("A5").Value = "Login"
$ws.Range("B5").Value = "Jason"
$ws.Range("C5").Value = "Implement credentials validation process"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1

$ws.Range("A6").Value = "Register"
$ws.Range("B6").Value = "JH"
$ws.Range("C6").Value = "Create register page"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0

$ws.Range("A7").Value = "Register"
$ws.Range("B7").Value = "JH"
$ws.Range("C7").Value = "Implement register process"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 3

$ws.Range("A8").Value = "Create home page"
$ws.Range("B8").Value = "JE"
$ws.Range("C8").Value = 'Create "Home Page"'
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0

$ws.Range("A9").Value = "Create home page"
$ws.Range("B9").Value = "JE"
$ws.Range("C9").Value = 'Implement "Home Page" functionality'
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 3

$ws.Range("A10").Value = "Create Inventory"
$ws.Range("B10").Value = "Vitor"
$ws.Range("C10").Value = 'Create "Inventory Page"'
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

$ws.Range("A11").Value = "Create Inventory"
$ws.Range("B11").Value = "Vitor"
$ws.Range("C11").Value = "Add user info at screen"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1

$ws.Range("A12").Value = "Create Inventory"
$ws.Range("B12").Value = "Vitor"
$ws.Range("C12").Value = "Implement add button for Raw Material"
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1

$ws.Range("A13").Value = "Create Inventory"
$ws.Range("B13").Value = "Vitor"
$ws.Range("C13").Value = "Implement edit button for Raw Material"
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1

$ws.Range("A14").Value = "Create Inventory"
$ws.Range("B14").Value = "Vitor"
$ws.Range("C14").Value = "Implement order material button for Raw Material"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 2

$ws.Range("A15").Value = "Create Inventory"
$ws.Range("B15").Value = "Vitor"
$ws.Range("C15").Value = "Implement remove button for Raw Material"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1

$ws.Range("A16").Value = "Create Inventory "
$ws.Range("B16").Value = "Vitor"
$ws.Range("C16").Value = "Implement Product recipe functionality"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1

$ws.Range("A17").Value = "Create Inventory"
$ws.Range("B17").Value = "Vitor"
$ws.Range("C17").Value = "Implement add button for Product"
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1

$ws.Range("A18").Value = "Create Inventory"
$ws.Range("B18").Value = "Vitor"
$ws.Range("C18").Value = "Implement edit button for Product"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1

$ws.Range("A19").Value = "Create Inventory"
$ws.Range("B19").Value = "Vitor"
$ws.Range("C19").Value = "Implement order material button for Product"
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1

$ws.Range("A20").Value = "Create Inventory"
$ws.Range("B20").Value = "Vitor"
$ws.Range("C20").Value = "Implement remove button for Product"
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1

# New row: Create Model Class for inventory
$ws.Range("A21").Value = "Create Inventory"
$ws.Range("B21").Value = "Vitor"
$ws.Range("C21").Value = "Create Model Class for inventory"
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 2

$ws.Range("A22").Value = "Create Orders"
$ws.Range("B22").Value = "JE"
$ws.Range("C22").Value = "Create Orders page"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0

# New row: Implement functionality for orders
$ws.Range("A23").Value = "Create Orders"
$ws.Range("B23").Value = "JE"
$ws.Range("C23").Value = "Implement functionality for orders"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0

# New row: Testing for orders functionality
$ws.Range("A24").Value = "Create Orders"
$ws.Range("B24").Value = "JE"
$ws.Range("C24").Value = "Testing for orders functionality"
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0

$ws.Range("A25").Value = "Create Orders"
$ws.Range("B25").Value = "JE"
$ws.Range("C25").Value = "Implement Open Orders list view"
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 2

$ws.Range("A26").Value = "Create Orders"
$ws.Range("B26").Value = "JE"
$ws.Range("C26").Value = "Implement Closed Orders list view"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 2

$ws.Range("A27").Value = "Create Orders"
$ws.Range("B27").ClearContents()
$ws.Range("C27").Value = "Create orders data base"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 3

# --- Column C width tweak (41.5 chars of XML column width) ---
$ws.Columns("C:C").ColumnWidth = 40.669221826364144

# --- Sheet view tweaks ---
$ws.Application.ActiveWindow.Zoom = 135
$ws.Range("C14").Select()
